# Generate Report for Handback
#
# The localization status workbook gets a "handback" pass recorded:
#  - Status text updates from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview + per-language sheets).
#  - Each per-language sheet (zh-cn, de-de) gets its "Latest Target File"
#    (F) / "Latest Handback File" (G) columns populated with hyperlinks
#    for both rows, mirroring the existing Source File Name (A) / Latest
#    Handoff File (D) hyperlinks.
#  - The "Latest Handback DateTime" (H) column is stamped with the
#    handback timestamp (per language).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrlTemplate  = "https://github.com/OpenLocalizationTest/oltest/blob/822ada9a4bcf95cb0233b57cd112bbaf790e1904/e2e/{0}.md"
$xlfUrlTemplate = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/{0}/ol-handoff/OpenLocalizationTestOrg/oltest.{1}/ci/ht/{2}"

$rows = @(
    @{ Id = "61b896cf-cc6b-4613-bae6-25589e9c641c"; Token = "0cb423db10d2ca3cac4e4e2e5696829bdf7b154d" },
    @{ Id = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6"; Token = "f46fd9bbdb5bce68e26b2f9491a78b463d29c64c" }
)

$languages = @(
    @{ Sheet = "zh-cn"; Commit = "aadc418136b2ab09f0d6912698843f0a83ab67e5"; HandbackTime = "2016-03-23 22:49:31" },
    @{ Sheet = "de-de"; Commit = "b2e6fc8ab0cd39017e653a934ebb2e062a927978"; HandbackTime = "2016-03-23 22:49:40" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $rowNum = 2

    foreach ($row in $rows) {
        $mdName  = $row.Id + ".md"
        $xlfName = $row.Id + "." + $row.Token + "." + $lang.Sheet + ".xlf"

        $mdUrl  = [string]::Format($mdUrlTemplate, $row.Id)
        $xlfUrl = [string]::Format($xlfUrlTemplate, $lang.Commit, $lang.Sheet, $xlfName)

        $statusCell = $ws.Cells.Item($rowNum, 3)   # column C - Status
        $statusCell.Value = $statusText

        $targetCell   = $ws.Cells.Item($rowNum, 6) # column F - Latest Target File
        $handbackCell = $ws.Cells.Item($rowNum, 7) # column G - Latest Handback File
        $dateCell     = $ws.Cells.Item($rowNum, 8) # column H - Latest Handback DateTime

        $ws.Hyperlinks.Add($targetCell, $mdUrl, "", "", $mdName) | Out-Null
        $ws.Hyperlinks.Add($handbackCell, $xlfUrl, "", "", $xlfName) | Out-Null

        $dateCell.Value = $lang.HandbackTime

        $rowNum = $rowNum + 1
    }
}

# Overview sheet mirrors the Status text in its per-language columns (B, C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText
